$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "A) Solutions:" -> "A) " + italic "(Pick one solution)"
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("A) Solutions:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'A) Solutions:' text"
}
# Re-create a plain Range from the found Start/End (InsertXML behaves as a
# true replace only on a freshly-constructed Range, not on Find's own
# Range object).
$target1 = $d.Range($find.Parent.Start, $find.Parent.End)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">A) </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>(Pick one solution)</w:t></w:r></w:p>'
$target1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: append new B)/C)/D) answer paragraphs after the final
# "b: H x 1" paragraph, moving the trailing line-break to the very end.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# The paragraph currently ends in "... x 1<break>" (break = vertical
# tab char right before the paragraph mark). Remove that break run so
# it can be re-created as its own paragraph after the new content.
$breakRange = $d.Range($lastRange.End - 2, $lastRange.End - 1)
$breakRange.Delete()

# Re-fetch the (now shorter) last paragraph and insert the new
# paragraphs right before its paragraph mark.
$lastPara2 = $d.Paragraphs.Last
$lastRange2 = $lastPara2.Range
$insertPoint = $d.Range($lastRange2.End - 1, $lastRange2.End - 1)

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>B) Matrix A: H x K</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">     Matrix B: K x D</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">     Bias </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>c :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> H</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">     Extra parameters per country: H x K + K x D + H.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">C) </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F04C"/></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">D) </w:t></w:r><w:r><w:t xml:space="preserve">Less </w:t></w:r><w:r><w:t xml:space="preserve">risk on </w:t></w:r><w:r><w:t>overfitting,</w:t></w:r><w:r><w:t xml:space="preserve"> less likely to memorize, training data, less memory consumption</w:t></w:r><w:r><w:t xml:space="preserve"> more generative</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>(pick any two of these reasons)</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:br/></w:r></w:p>'
$insertPoint.InsertXML($xml2)

Write-Output "edit complete"
